$d = $word.ActiveDocument

# Paragraph 40 is the last paragraph: "Uppd. Designdoc: ... en plattform"
# It currently holds the _GoBack bookmark around its single run; remove that
# bookmark here -- it will be re-added further down, around the new final
# text ("... halvsekund till").
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$p40 = $d.Paragraphs.Item(40)
$full = $d.Range($p40.Range.Start, $p40.Range.End - 1)
$full.Text = ""

$ins = $d.Range($p40.Range.Start, $p40.Range.Start)
$ins.InsertAfter("Uppd. Designdoc: Lägg till spec. för hur långt ")
$ins.Font.Bold = 1
$ins2 = $d.Range($ins.End, $ins.End)
$ins2.InsertAfter("en del kan gå ut på en plattform")
$ins2.Font.Bold = 1

# New paragraph: Robotar i bakgrunden
$p40 = $d.Paragraphs.Item(40)
$p40.Range.InsertParagraphAfter()
$p41 = $d.Paragraphs.Item(41)
$p41.Range.InsertBefore("Robotar i bakgrunden")

# New paragraph: Tid börjar när Stix börjar röra på sig
$p41.Range.InsertParagraphAfter()
$p42 = $d.Paragraphs.Item(42)
$p42.Range.InsertBefore("Tid börjar när Stix börjar röra på sig")

# New paragraph: Paus efter avklarad bana då fanfar-musik spelas
$p42.Range.InsertParagraphAfter()
$p43 = $d.Paragraphs.Item(43)
$p43.Range.InsertBefore("Paus efter avklarad bana då fanfar-musik spelas")

# New paragraph: Den svagaste tilen ... halvsekund [ till]
$p43.Range.InsertParagraphAfter()
$p44 = $d.Paragraphs.Item(44)
$p44.Range.InsertBefore("Den svagaste tilen skakar/smular sönder och kan därför fortfarande användas i ngn halvsekund")

$p44 = $d.Paragraphs.Item(44)
$tail = $d.Range($p44.Range.End - 1, $p44.Range.End - 1)
$tail.InsertAfter(" till")
$tail.Font.Bold = 1

# Re-add the _GoBack bookmark, now collapsed at the end of paragraph 44's text
$p44 = $d.Paragraphs.Item(44)
$gobackPos = $p44.Range.End - 1
$gobackRange = $d.Range($gobackPos, $gobackPos)
$d.Bookmarks.Add("_GoBack", $gobackRange)

# Final new empty paragraph after paragraph 44
$p44 = $d.Paragraphs.Item(44)
$p44.Range.InsertParagraphAfter()
